$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "242.90"
Set-TextValue "G2" "14"
Set-TextValue "D3" "23.01"
Set-TextValue "G3" "14"
Set-TextValue "D4" "5.400"
Set-TextValue "G4" "14"
Set-TextValue "D5" "0.05901"
Set-TextValue "G5" "14"
Set-TextValue "G6" "14"
Set-TextValue "D7" "6.589"
Set-TextValue "G7" "14"
Set-TextValue "D8" "0.8104"
Set-TextValue "G8" "14"
Set-TextValue "D9" "0.9190"
Set-TextValue "G9" "14"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1422"
Set-TextValue "E10" "9WazirXWRX"
Set-TextValue "G10" "14"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07399"
Set-TextValue "E11" "10MandalaExchangeTokenMDX"
Set-TextValue "G11" "14"
Set-TextValue "B12" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03272"
Set-TextValue "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "G12" "14"
Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03060"
Set-TextValue "E13" "12BitrueCoinBTR"
Set-TextValue "G13" "14"
Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09337"
Set-TextValue "E14" "13BitMartTokenBMX"
Set-TextValue "G14" "14"
Set-TextValue "B15" "MCDex"
Set-TextValue "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.849"
Set-TextValue "E15" "14MCDexMCB"
Set-TextValue "G15" "14"
Set-TextValue "B16" "BitForexToken"
Set-TextValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001558"
Set-TextValue "E16" "15BitForexTokenBF"
Set-TextValue "G16" "14"
Set-TextValue "B17" "CoinExToken"
Set-TextValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04663"
Set-TextValue "E17" "16CoinExTokenCET"
Set-TextValue "G17" "14"
Set-TextValue "B18" "One"
Set-TextValue "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005953"
Set-TextValue "E18" "17OneONE"
Set-TextValue "G18" "14"
Set-TextValue "D19" "0.005895"
Set-TextValue "G19" "14"
Set-TextValue "D20" "0.001277"
Set-TextValue "E20" "19BitKanKANBestin24h"
Set-TextValue "G20" "14"
Set-TextValue "D21" "0.004894"
Set-TextValue "G21" "14"
Set-TextValue "D22" "0.00009502"
Set-TextValue "G22" "14"
Set-TextValue "D23" "3.614"
Set-TextValue "G23" "14"
Set-TextValue "D24" "2.144"
Set-TextValue "G24" "14"
Set-TextValue "D25" "0.3230"
Set-TextValue "G25" "14"
Set-TextValue "G26" "14"
Set-TextValue "G27" "14"
Set-TextValue "G28" "14"
Set-TextValue "G29" "14"
Set-TextValue "G30" "14"
Set-TextValue "G31" "14"
Set-TextValue "G32" "14"
Set-TextValue "G33" "14"
Set-TextValue "G34" "14"
Set-TextValue "G35" "14"
Set-TextValue "G36" "14"
Set-TextValue "G37" "14"
Set-TextValue "G38" "14"
Set-TextValue "G39" "14"
Set-TextValue "D40" "0.03960"
Set-TextValue "G40" "14"
Set-TextValue "D41" "0.006177"
Set-TextValue "G41" "14"
Set-TextValue "G42" "14"
Set-TextValue "G43" "14"
Set-TextValue "D44" "0.008110"
Set-TextValue "G44" "14"
Set-TextValue "D45" "0.00005194"
Set-TextValue "G45" "14"
Set-TextValue "G46" "14"
Set-TextValue "D47" "0.7504"
Set-TextValue "G47" "14"
Set-TextValue "D48" "0.002277"
Set-TextValue "G48" "14"
Set-TextValue "D49" "0.00002100"
Set-TextValue "G49" "14"
Set-TextValue "D50" "0.0002000"
Set-TextValue "G50" "14"
Set-TextValue "G51" "14"

Write-Host "Applied 109 cell updates"